$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8541.5
$ws.Range("I2").Value = 507.42856
$ws.Range("K2").Value = 507.42856
$ws.Range("M2").Value = -394.42856

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6283.0435
$ws.Range("I11").Value = 6283.0435
$ws.Range("K11").Value = 6283.0435
$ws.Range("M11").Value = -6143.0435

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 9892.941000000001
$ws.Range("J28").Value = 17833
$ws.Range("L28").Value = 17833
$ws.Range("N28").Value = -18803

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3620.6667
$ws.Range("I29").Value = 1931
$ws.Range("K29").Value = 5793
$ws.Range("M29").Value = -5512

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 7377
$ws.Range("I98").Value = 7377
$ws.Range("K98").Value = 7377
$ws.Range("M98").Value = -5879

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 7377
$ws.Range("I122").Value = 7377
$ws.Range("K122").Value = 22131
$ws.Range("M122").Value = -19681

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1281.6875
$ws.Range("I129").Value = 916.38464
$ws.Range("J129").Value = 2864.6667
$ws.Range("K129").Value = 2749.15392
$ws.Range("L129").Value = 8594.000100000001
$ws.Range("M129").Value = 2250.84608
$ws.Range("N129").Value = -18594.0001

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2397.1667
$ws.Range("I131").Value = 1386.6
$ws.Range("J131").Value = 7450
$ws.Range("K131").Value = 4159.799999999999
$ws.Range("L131").Value = 22350
$ws.Range("M131").Value = 880.2000000000007
$ws.Range("N131").Value = -32430

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5630.92
$ws.Range("I137").Value = 936.17645
$ws.Range("K137").Value = 2808.52935
$ws.Range("M137").Value = -258.5293500000002

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1001
$ws.Range("I4").Value = 1001
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -885
$ws.Range("N4").ClearContents()

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3356.4285
$ws.Range("I63").Value = 4299.3335
$ws.Range("J63").Value = 2649.25
$ws.Range("K63").Value = 4299.3335
$ws.Range("L63").Value = 2649.25
$ws.Range("M63").Value = -3613.3335
$ws.Range("N63").Value = -4021.25

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3356.4285
$ws.Range("I66").Value = 4299.3335
$ws.Range("J66").Value = 2649.25
$ws.Range("K66").Value = 21496.6675
$ws.Range("L66").Value = 13246.25
$ws.Range("M66").Value = -18064.6675
$ws.Range("N66").Value = -20110.25

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1612935.2
$ws.Range("I97").Value = 1612935.2
$ws.Range("K97").Value = 1612935.2
$ws.Range("M97").Value = -1612439.2

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1312.1333
$ws.Range("I122").Value = 1312.1333
$ws.Range("K122").Value = 3936.3999
$ws.Range("M122").Value = -1486.3999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3805.375
$ws.Range("J132").Value = 5727.5713
$ws.Range("L132").Value = 17182.7139
$ws.Range("N132").Value = -22242.7139

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 98653.91
$ws.Range("I86").Value = 4930.625
$ws.Range("K86").Value = 4930.625
$ws.Range("M86").Value = -3807.625

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 98653.91
$ws.Range("I89").Value = 4930.625
$ws.Range("K89").Value = 24653.125
$ws.Range("M89").Value = -19037.125

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7401.3
$ws.Range("I134").Value = 2633.3333
$ws.Range("K134").Value = 7899.999899999999
$ws.Range("M134").Value = -5364.999899999999

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 337.04166
$ws.Range("I22").Value = 258.83334
$ws.Range("K22").Value = 258.83334
$ws.Range("M22").Value = 91.16665999999998

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4511.727
$ws.Range("J31").Value = 5631.2666
$ws.Range("L31").Value = 5631.2666
$ws.Range("N31").Value = -6221.2666

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4511.727
$ws.Range("J34").Value = 5631.2666
$ws.Range("L34").Value = 5631.2666
$ws.Range("N34").Value = -6035.2666

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5059.1577
$ws.Range("I58").Value = 3827.7778
$ws.Range("J58").Value = 6167.4
$ws.Range("K58").Value = 3827.7778
$ws.Range("L58").Value = 6167.4
$ws.Range("M58").Value = -3624.7778
$ws.Range("N58").Value = -6573.4

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 56245.22
$ws.Range("I132").Value = 2363.182
$ws.Range("K132").Value = 7089.545999999999
$ws.Range("M132").Value = -4559.545999999999

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5599.2
$ws.Range("I134").Value = 5325.484
$ws.Range("K134").Value = 15976.452
$ws.Range("M134").Value = -13441.452

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5059.1577
$ws.Range("I136").Value = 3827.7778
$ws.Range("J136").Value = 6167.4
$ws.Range("K136").Value = 11483.3334
$ws.Range("L136").Value = 18502.2
$ws.Range("M136").Value = -8933.3334
$ws.Range("N136").Value = -23602.2

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3443
$ws.Range("I137").Value = 1610
$ws.Range("J137").Value = 4228.5713
$ws.Range("K137").Value = 4830
$ws.Range("L137").Value = 12685.7139
$ws.Range("M137").Value = 270
$ws.Range("N137").Value = -22885.7139

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3850353.2
$ws.Range("J138").Value = 6216.375
$ws.Range("L138").Value = 18649.125
$ws.Range("N138").Value = -28929.125

# GSM row 55
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 44000
$ws.Range("J55").Value = 44000
$ws.Range("L55").Value = 44000
$ws.Range("N55").Value = -44654

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2393.6875
$ws.Range("I126").Value = 2374.6924
$ws.Range("J126").Value = 2476
$ws.Range("K126").Value = 7124.0772
$ws.Range("L126").Value = 7428
$ws.Range("M126").Value = -4654.0772
$ws.Range("N126").Value = -12368

# LTW row 76
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 11329.667
$ws.Range("J76").Value = 11329.667
$ws.Range("L76").Value = 11329.667
$ws.Range("N76").Value = -12005.667

# LTW row 79
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 11329.667
$ws.Range("J79").Value = 11329.667
$ws.Range("L79").Value = 11329.667
$ws.Range("N79").Value = -13669.667

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6899
$ws.Range("I93").Value = 6427.143
$ws.Range("J93").Value = 8000
$ws.Range("K93").Value = 6427.143
$ws.Range("L93").Value = 8000
$ws.Range("M93").Value = -5179.143
$ws.Range("N93").Value = -10496

# LTW row 101
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 21000
$ws.Range("J101").Value = 21000
$ws.Range("L101").Value = 21000
$ws.Range("N101").Value = -27490

# WVR row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41144

# WVR row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 10979.6
$ws.Range("I52").Value = 4285.2856
$ws.Range("J52").Value = 26599.666
$ws.Range("K52").Value = 4285.2856
$ws.Range("L52").Value = 26599.666
$ws.Range("M52").Value = -4059.2856
$ws.Range("N52").Value = -27051.666

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4959
$ws.Range("I136").Value = 2757.9333
$ws.Range("K136").Value = 8273.7999
$ws.Range("M136").Value = -5723.7999
